# Apply the updated yearly financial figures to the CCIH worksheet.
# Mirrors the author's "Doing Updates for Financials" data refresh:
# every historical revenue/expense/cash-flow figure on the sheet is
# nudged to its newer reported value (cell formatting/layout is left
# untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Cell = "D8"; Value = 126500 },
    @{ Cell = "E8"; Value = 156500 },
    @{ Cell = "F8"; Value = 200900 },
    @{ Cell = "G8"; Value = 205400 },
    @{ Cell = "H8"; Value = 163700 },
    @{ Cell = "I8"; Value = 120800 },
    @{ Cell = "J8"; Value = 91800 },
    @{ Cell = "D9"; Value = 116000 },
    @{ Cell = "E9"; Value = 160000 },
    @{ Cell = "F9"; Value = 154600 },
    @{ Cell = "G9"; Value = 143400 },
    @{ Cell = "H9"; Value = 112300 },
    @{ Cell = "I9"; Value = 83000 },
    @{ Cell = "J9"; Value = 64200 },
    @{ Cell = "D10"; Value = 10500 },
    @{ Cell = "E10"; Value = -3500 },
    @{ Cell = "F10"; Value = 46300 },
    @{ Cell = "G10"; Value = 62000 },
    @{ Cell = "H10"; Value = 51400 },
    @{ Cell = "I10"; Value = 37700 },
    @{ Cell = "J10"; Value = 27600 },
    @{ Cell = "D12"; Value = 12100 },
    @{ Cell = "E12"; Value = 15400 },
    @{ Cell = "F12"; Value = 13400 },
    @{ Cell = "G12"; Value = 15900 },
    @{ Cell = "H12"; Value = 14000 },
    @{ Cell = "I12"; Value = 9300 },
    @{ Cell = "J12"; Value = 6100 },
    @{ Cell = "D14"; Value = 3800 },
    @{ Cell = "E14"; Value = 61900 },
    @{ Cell = "J14"; Value = 1100 },
    @{ Cell = "F15"; Value = 3300 },
    @{ Cell = "D17"; Value = 167800 },
    @{ Cell = "E17"; Value = 293400 },
    @{ Cell = "F17"; Value = 218500 },
    @{ Cell = "G17"; Value = 208000 },
    @{ Cell = "H17"; Value = 169000 },
    @{ Cell = "I17"; Value = 121900 },
    @{ Cell = "J17"; Value = 90600 },
    @{ Cell = "D18"; Value = -41200 },
    @{ Cell = "E18"; Value = -136900 },
    @{ Cell = "F18"; Value = -17700 },
    @{ Cell = "G18"; Value = -2600 },
    @{ Cell = "H18"; Value = -5300 },
    @{ Cell = "E20"; Value = 3600 },
    @{ Cell = "F20"; Value = 3100 },
    @{ Cell = "D21"; Value = -41800 },
    @{ Cell = "E21"; Value = -109700 },
    @{ Cell = "F21"; Value = 8900 },
    @{ Cell = "G21"; Value = 14200 },
    @{ Cell = "H21"; Value = 4800 },
    @{ Cell = "I21"; Value = 7500 },
    @{ Cell = "J21"; Value = 10200 },
    @{ Cell = "D22"; Value = 2800 },
    @{ Cell = "F22"; Value = 2000 },
    @{ Cell = "D23"; Value = -46200 },
    @{ Cell = "E23"; Value = -135100 },
    @{ Cell = "F23"; Value = -16500 },
    @{ Cell = "G23"; Value = -1500 },
    @{ Cell = "H23"; Value = -4900 },
    @{ Cell = "D24"; Value = 8900 },
    @{ Cell = "F24"; Value = -3400 },
    @{ Cell = "J24"; Value = 1700 },
    @{ Cell = "D26"; Value = -55100 },
    @{ Cell = "E26"; Value = -135700 },
    @{ Cell = "F26"; Value = -13200 },
    @{ Cell = "H26"; Value = -5100 },
    @{ Cell = "D27"; Value = -54800 },
    @{ Cell = "E27"; Value = -135600 },
    @{ Cell = "F27"; Value = -13200 },
    @{ Cell = "H27"; Value = -5100 },
    @{ Cell = "E32"; Value = -3600 },
    @{ Cell = "F32"; Value = -3100 },
    @{ Cell = "D33"; Value = -54800 },
    @{ Cell = "E33"; Value = -135600 },
    @{ Cell = "F33"; Value = -13200 },
    @{ Cell = "H33"; Value = -5100 },
    @{ Cell = "D35"; Value = -54800 },
    @{ Cell = "E35"; Value = -135600 },
    @{ Cell = "F35"; Value = -13200 },
    @{ Cell = "H35"; Value = -5100 },
    @{ Cell = "D41"; Value = 15800 },
    @{ Cell = "E41"; Value = 20000 },
    @{ Cell = "F41"; Value = 90100 },
    @{ Cell = "G41"; Value = 55800 },
    @{ Cell = "H41"; Value = 50200 },
    @{ Cell = "I41"; Value = 47100 },
    @{ Cell = "J41"; Value = 58300 },
    @{ Cell = "F42"; Value = 3900 },
    @{ Cell = "H42"; Value = 3700 },
    @{ Cell = "D43"; Value = 24200 },
    @{ Cell = "G43"; Value = 47400 },
    @{ Cell = "H43"; Value = 46200 },
    @{ Cell = "I43"; Value = 35400 },
    @{ Cell = "J43"; Value = 23200 },
    @{ Cell = "D45"; Value = 117900 },
    @{ Cell = "E45"; Value = 201300 },
    @{ Cell = "F45"; Value = 163900 },
    @{ Cell = "G45"; Value = 21400 },
    @{ Cell = "H45"; Value = 16800 },
    @{ Cell = "I45"; Value = 6700 },
    @{ Cell = "D46"; Value = 158000 },
    @{ Cell = "E46"; Value = 257500 },
    @{ Cell = "F46"; Value = 294900 },
    @{ Cell = "G46"; Value = 128300 },
    @{ Cell = "H46"; Value = 116800 },
    @{ Cell = "I46"; Value = 89200 },
    @{ Cell = "J46"; Value = 85800 },
    @{ Cell = "D47"; Value = 4500 },
    @{ Cell = "E47"; Value = 5100 },
    @{ Cell = "F47"; Value = 7400 },
    @{ Cell = "G47"; Value = 7000 },
    @{ Cell = "H47"; Value = 5000 },
    @{ Cell = "I47"; Value = 14500 },
    @{ Cell = "J47"; Value = 16100 },
    @{ Cell = "D48"; Value = 7900 },
    @{ Cell = "F48"; Value = 74200 },
    @{ Cell = "G48"; Value = 62200 },
    @{ Cell = "H48"; Value = 35700 },
    @{ Cell = "I48"; Value = 26600 },
    @{ Cell = "J48"; Value = 23000 },
    @{ Cell = "D49"; Value = 4900 },
    @{ Cell = "G49"; Value = 8900 },
    @{ Cell = "H49"; Value = 8400 },
    @{ Cell = "D52"; Value = 63100 },
    @{ Cell = "E52"; Value = 5400 },
    @{ Cell = "F52"; Value = 10500 },
    @{ Cell = "G52"; Value = 50500 },
    @{ Cell = "H52"; Value = 7400 },
    @{ Cell = "I52"; Value = 3000 },
    @{ Cell = "D54"; Value = 238400 },
    @{ Cell = "E54"; Value = 268000 },
    @{ Cell = "F54"; Value = 388600 },
    @{ Cell = "G54"; Value = 256900 },
    @{ Cell = "H54"; Value = 173200 },
    @{ Cell = "I54"; Value = 133700 },
    @{ Cell = "J54"; Value = 125700 },
    @{ Cell = "D57"; Value = 54600 },
    @{ Cell = "E57"; Value = 44800 },
    @{ Cell = "F57"; Value = 30500 },
    @{ Cell = "G57"; Value = 38000 },
    @{ Cell = "H57"; Value = 30400 },
    @{ Cell = "I57"; Value = 15900 },
    @{ Cell = "J57"; Value = 8900 },
    @{ Cell = "D58"; Value = 12700 },
    @{ Cell = "E58"; Value = 15700 },
    @{ Cell = "F58"; Value = 11500 },
    @{ Cell = "G58"; Value = 12000 },
    @{ Cell = "H58"; Value = 8900 },
    @{ Cell = "D59"; Value = 212800 },
    @{ Cell = "E59"; Value = 220500 },
    @{ Cell = "F59"; Value = 222100 },
    @{ Cell = "G59"; Value = 78300 },
    @{ Cell = "H59"; Value = 36700 },
    @{ Cell = "I59"; Value = 19000 },
    @{ Cell = "J59"; Value = 17000 },
    @{ Cell = "D60"; Value = 280100 },
    @{ Cell = "E60"; Value = 281000 },
    @{ Cell = "F60"; Value = 264100 },
    @{ Cell = "G60"; Value = 128200 },
    @{ Cell = "H60"; Value = 76000 },
    @{ Cell = "I60"; Value = 35000 },
    @{ Cell = "J60"; Value = 26600 },
    @{ Cell = "D61"; Value = 31600 },
    @{ Cell = "E61"; Value = 6500 },
    @{ Cell = "F61"; Value = 16100 },
    @{ Cell = "G61"; Value = 4800 },
    @{ Cell = "E62"; Value = 1700 },
    @{ Cell = "F62"; Value = 1300 },
    @{ Cell = "D66"; Value = 312300 },
    @{ Cell = "E66"; Value = 289100 },
    @{ Cell = "F66"; Value = 281600 },
    @{ Cell = "G66"; Value = 133000 },
    @{ Cell = "H66"; Value = 76300 },
    @{ Cell = "I66"; Value = 35500 },
    @{ Cell = "J66"; Value = 26600 },
    @{ Cell = "D72"; Value = -307900 },
    @{ Cell = "E72"; Value = -249500 },
    @{ Cell = "F72"; Value = -98300 },
    @{ Cell = "G72"; Value = -81500 },
    @{ Cell = "H72"; Value = -77600 },
    @{ Cell = "I72"; Value = -72500 },
    @{ Cell = "J72"; Value = -70000 },
    @{ Cell = "D76"; Value = -74000 },
    @{ Cell = "E76"; Value = -21100 },
    @{ Cell = "F76"; Value = 107000 },
    @{ Cell = "G76"; Value = 123900 },
    @{ Cell = "H76"; Value = 96900 },
    @{ Cell = "I76"; Value = 98200 },
    @{ Cell = "J76"; Value = 99100 },
    @{ Cell = "D81"; Value = -54800 },
    @{ Cell = "E81"; Value = -135600 },
    @{ Cell = "F81"; Value = -13200 },
    @{ Cell = "H81"; Value = -5100 },
    @{ Cell = "E83"; Value = 23600 },
    @{ Cell = "F83"; Value = 23400 },
    @{ Cell = "G83"; Value = 14400 },
    @{ Cell = "H83"; Value = 9200 },
    @{ Cell = "I83"; Value = 8900 },
    @{ Cell = "J83"; Value = 10100 },
    @{ Cell = "D89"; Value = -14700 },
    @{ Cell = "E89"; Value = -27800 },
    @{ Cell = "F89"; Value = 114300 },
    @{ Cell = "G89"; Value = 41600 },
    @{ Cell = "H89"; Value = 18000 },
    @{ Cell = "I89"; Value = 4700 },
    @{ Cell = "J89"; Value = 8700 },
    @{ Cell = "D91"; Value = -2300 },
    @{ Cell = "E91"; Value = -8800 },
    @{ Cell = "F91"; Value = -13400 },
    @{ Cell = "G91"; Value = -37600 },
    @{ Cell = "H91"; Value = -14500 },
    @{ Cell = "I91"; Value = -15200 },
    @{ Cell = "J91"; Value = -11500 },
    @{ Cell = "D94"; Value = -13300 },
    @{ Cell = "E94"; Value = -30000 },
    @{ Cell = "F94"; Value = -64200 },
    @{ Cell = "G94"; Value = -62600 },
    @{ Cell = "H94"; Value = -16800 },
    @{ Cell = "I94"; Value = -14800 },
    @{ Cell = "J94"; Value = -28600 },
    @{ Cell = "D100"; Value = 22100 },
    @{ Cell = "E100"; Value = -12600 },
    @{ Cell = "F100"; Value = -15800 },
    @{ Cell = "G100"; Value = 26400 },
    @{ Cell = "J100"; Value = -8800 },
    @{ Cell = "D101"; Value = -1600 },
    @{ Cell = "E101"; Value = 2200 },
    @{ Cell = "D102"; Value = -7400 },
    @{ Cell = "E102"; Value = -68200 },
    @{ Cell = "F102"; Value = 35700 },
    @{ Cell = "G102"; Value = 5600 },
    @{ Cell = "H102"; Value = 3100 },
    @{ Cell = "I102"; Value = -11200 },
    @{ Cell = "J102"; Value = -29700 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
